$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new experiment row (row 17) ---------------------------------
# Order matters: new shared-strings get appended in the order the values are
# first assigned, so write B17, then A17, then O15/O16/O17 to reproduce the
# same shared-string ordering as the target workbook.

$ws.Range("B17").Value = "PPO use step distance reward + multiply critic lr + train every episode + summed hit wall penalty + actor 2 layers + 3 frames vs. Random"

# A17 needs the same "short name" cell style used by A16 (style index 6,
# grey fill). Copy formatting only (not the value) from A16, then overwrite
# the value so the existing style entry is reused instead of a new one
# being created.
[void]$ws.Range("A16").Copy()
[void]$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "SD-SHW-3F"

$ws.Range("O15").Value = "map*_use_step_dist_normed"
$ws.Range("O16").Value = "map*_use_step_dist_summed_hit_wall"
$ws.Range("O17").Value = "map*_use_step_dist_summed_hit_wall_[actor]2layers_[frames]3"

# --- Widen column O so the long "Run dir" strings are readable -----------
$ws.Columns.Item(15).ColumnWidth = 60.428571428571427

# --- Update the view: scroll to B1, zoom 87%, select E19 -----------------
[void]$ws.Range("B1").Select()
$excel.ActiveWindow.Zoom = 87
[void]$ws.Range("E19").Select()
